$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Meter Consumption ID" column (E) and the "Project Name" column (B)
# -- delete the right-most one first so the left column letter stays valid.
$ws.Columns("E:E").Delete()
$ws.Columns("B:B").Delete()

# The Start/End Date values (now columns D:E) pick up a custom date+time
# display format instead of the plain short-date format they used to have.
$ws.Range("D2:E2").NumberFormat = "m/d/yy\ h:mm;@"

# Page was set to print in portrait orientation.
$ws.PageSetup.Orientation = 1

# Restore the selection that was active when the workbook was saved.
$ws.Range("D6").Select()
